# Update the division problems in the table to the new set of values.
# Each old value is unique within the document, so a simple Find/Replace
# (one occurrence at a time) for each pair is safe and order-independent.

$d = $word.ActiveDocument

$pairs = @(
    @{old = "937÷9="; new = "641÷3="},
    @{old = "272÷5="; new = "642÷4="},
    @{old = "598÷8="; new = "682÷8="},
    @{old = "432÷6="; new = "391÷5="},
    @{old = "342÷9="; new = "916÷3="},
    @{old = "707÷4="; new = "772÷8="},
    @{old = "253÷9="; new = "297÷8="},
    @{old = "681÷7="; new = "136÷6="},
    @{old = "845÷3="; new = "353÷2="},
    @{old = "116÷6="; new = "564÷2="},
    @{old = "378÷8="; new = "231÷2="},
    @{old = "427÷3="; new = "370÷3="},
    @{old = "646÷8="; new = "744÷7="},
    @{old = "332÷9="; new = "885÷7="},
    @{old = "448÷4="; new = "521÷8="},
    @{old = "129÷7="; new = "948÷2="},
    @{old = "374÷2="; new = "268÷8="},
    @{old = "148÷3="; new = "845÷3="},
    @{old = "100÷2="; new = "106÷7="},
    @{old = "940÷5="; new = "246÷7="},
    @{old = "241÷5="; new = "887÷5="},
    @{old = "951÷5="; new = "716÷9="},
    @{old = "454÷4="; new = "590÷3="},
    @{old = "194÷7="; new = "478÷4="},
    @{old = "239÷3="; new = "650÷7="}
)

foreach ($pair in $pairs) {
    $d.Content.Find.Execute($pair.old, $true, $true, $false, $false, $false,
                             $true, 1, $false, $pair.new, 2)
}
